# Add the new "07/02/2022 - 13/02/2022" week row (row 15) to every
# per-level sheet, plus the aggregated "Totale casi" sheet.

$wb = $excel.ActiveWorkbook

$newWeek = "07/02/2022 - 13/02/2022"

# Infanzia: Personale scolastico=27, Alunni=173, Totale=200
$ws = $wb.Worksheets.Item("Infanzia")
$ws.Range("A15").Value = $newWeek
$ws.Range("B15").Value = 27
$ws.Range("C15").Value = 173
$ws.Range("D15").Value = 200
$ws.Cells.Item(16, 2).Select()

# Primaria: Personale scolastico=12, Alunni=226, Totale=238
$ws = $wb.Worksheets.Item("Primaria")
$ws.Range("A15").Value = $newWeek
$ws.Range("B15").Value = 12
$ws.Range("C15").Value = 226
$ws.Range("D15").Value = 238
$ws.Cells.Item(15, 5).Select()

# Media: Personale scolastico=13, Alunni=183, Totale=196
$ws = $wb.Worksheets.Item("Media")
$ws.Range("A15").Value = $newWeek
$ws.Range("B15").Value = 13
$ws.Range("C15").Value = 183
$ws.Range("D15").Value = 196
$ws.Cells.Item(15, 5).Select()

# Superiore: Personale scolastico=14, Alunni=239, Totale=253
$ws = $wb.Worksheets.Item("Superiore")
$ws.Range("A15").Value = $newWeek
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 239
$ws.Range("D15").Value = 253
$ws.Cells.Item(16, 2).Select()

# Totale casi: Personale scolastico=66, Alunni=821, Totale=887
$ws = $wb.Worksheets.Item("Totale casi")
$ws.Range("A15").Value = $newWeek
$ws.Range("B15").Value = 66
$ws.Range("C15").Value = 821
$ws.Range("D15").Value = 887
$ws.Cells.Item(16, 2).Select()

# Keep "Totale casi" as the active/selected tab, matching the source file.
$ws.Activate()
